$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G12").Value = "fsdf"
$ws.Range("E6").Value = "dfdf"

$ws.Range("E6").Select()
